$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("6:11").Insert()

$ws.Range("A6").Value = 9
$ws.Range("B6").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C6").Value = "Metropolitana"
$ws.Range("D6").Value = 44490
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107002
$ws.Range("J6").Value = "Chirimoya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Cuarta"
$ws.Range("M6").Value = 330
$ws.Range("N6").Value = 1200
$ws.Range("O6").Value = 1200
$ws.Range("P6").Value = 1200
$ws.Range("Q6").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 1200
$ws.Range("T6").Value = 1

$ws.Range("A7").Value = 9
$ws.Range("B7").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44490
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107002
$ws.Range("J7").Value = "Chirimoya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 440
$ws.Range("N7").Value = 17600
$ws.Range("O7").Value = 17600
$ws.Range("P7").Value = 17600
$ws.Range("Q7").Value = "`$/bandeja 8 kilos"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 2200
$ws.Range("T7").Value = 8

$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44490
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107002
$ws.Range("J8").Value = "Chirimoya"
$ws.Range("K8").Value = "Cultivar IV Región"
$ws.Range("L8").Value = "Extra (doble especial)"
$ws.Range("M8").Value = 410
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 24000
$ws.Range("Q8").Value = "`$/bandeja 8 kilos"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 3000
$ws.Range("T8").Value = 8

$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44490
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107002
$ws.Range("J9").Value = "Chirimoya"
$ws.Range("K9").Value = "Cultivar IV Región"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 350
$ws.Range("N9").Value = 16000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 16000
$ws.Range("Q9").Value = "`$/bandeja 8 kilos"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 2000
$ws.Range("T9").Value = 8

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44490
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = "Otros"
$ws.Range("I10").Value = 100107002
$ws.Range("J10").Value = "Chirimoya"
$ws.Range("K10").Value = "Cultivar IV Región"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 14400
$ws.Range("O10").Value = 14400
$ws.Range("P10").Value = 14400
$ws.Range("Q10").Value = "`$/bandeja 8 kilos"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 1800
$ws.Range("T10").Value = 8

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44490
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100107
$ws.Range("H11").Value = "Otros"
$ws.Range("I11").Value = 100107002
$ws.Range("J11").Value = "Chirimoya"
$ws.Range("K11").Value = "Cultivar IV Región"
$ws.Range("L11").Value = "Tercera"
$ws.Range("M11").Value = 290
$ws.Range("N11").Value = 1400
$ws.Range("O11").Value = 1400
$ws.Range("P11").Value = 1400
$ws.Range("Q11").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 1400
$ws.Range("T11").Value = 1
